# Auto-generated edit script: applies the "Updated cryptos list" diff
# (price / volume refresh + two coin-row reorderings) cell by cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) hold numeric-looking text that must
# stay plain text (inline/shared string), matching the original cells which
# have no explicit number format. Forcing NumberFormat to "@" (Text) before
# assigning the value stops Excel from auto-converting strings such as
# "645.01" or "  -5.55%  " into real numbers; ClearFormats() afterwards drops
# the now-unneeded Text format so the cell keeps its original (default) style.
function Set-TextValue($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue "D2" "69.331.32"
Set-TextValue "E2" "  -0.04%  "
Set-TextValue "D3" "3.667.26"
Set-TextValue "E3" "  -0.62%  "
Set-TextValue "E4" "  +0.14%  "
Set-TextValue "D5" "645.01"
Set-TextValue "E5" "  -5.55%  "
Set-TextValue "D6" "158.96"
Set-TextValue "E6" "  -0.45%  "
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "E8" "  +0.09%  "
Set-TextValue "E9" "  -1.03%  "
Set-TextValue "D10" "7.10"
Set-TextValue "E10" "  -0.32%  "
Set-TextValue "D11" "0.440"
Set-TextValue "E11" "  +0.17%  "
Set-TextValue "E12" "  -0.27%  "
Set-TextValue "D13" "4.291.83"
Set-TextValue "E13" "  -0.45%  "
Set-TextValue "D14" "32.50"
Set-TextValue "E14" "  -0.01%  "
Set-TextValue "D15" "3.668.20"
Set-TextValue "E15" "  -0.45%  "
Set-TextValue "D16" "69.383.48"
Set-TextValue "E16" "  +0.05%  "
Set-TextValue "D18" "15.99"
Set-TextValue "E18" "  -0.70%  "
Set-TextValue "D19" "6.47"
Set-TextValue "E19" "  -0.25%  "
Set-TextValue "D20" "465.67"
Set-TextValue "E20" "  -0.89%  "
Set-TextValue "D21" "9.83"
Set-TextValue "E21" "  -1.06%  "
Set-TextValue "E22" "  -1.81%  "
Set-TextValue "D23" "79.40"
Set-TextValue "E23" "  -0.81%  "
Set-TextValue "D24" "3.818.82"
Set-TextValue "E24" "  -0.45%  "
Set-TextValue "E25" "  -0.03%  "
Set-TextValue "E26" "  +0.46%  "
Set-TextValue "D27" "10.84"
Set-TextValue "E27" "  -1.03%  "
Set-TextValue "D28" "8.93"
Set-TextValue "E28" "  -2.55%  "
Set-TextValue "E29" "  -2.94%  "
Set-TextValue "D30" "1.69"
Set-TextValue "E30" "  -3.12%  "
Set-TextValue "D31" "2.00"
Set-TextValue "E31" "  -0.06%  "
Set-TextValue "D32" "1.00"
Set-TextValue "E32" "  +0.06%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "26.62"
Set-TextValue "E33" "  -1.24%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D34" "6.46"
Set-TextValue "E34" "  -2.81%  "
Set-TextValue "D35" "0.162"
Set-TextValue "E35" "  +3.94%  "
Set-TextValue "D36" "3.661.95"
Set-TextValue "E36" "  -0.42%  "
Set-TextValue "D37" "8.40"
Set-TextValue "E37" "  +1.44%  "
Set-TextValue "D39" "5.90"
Set-TextValue "E39" "  -6.31%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D40" "177.92"
Set-TextValue "E40" "  +4.68%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D41" "1.00"
Set-TextValue "E41" "  +0.08%  "
Set-TextValue "D42" "0.0897"
Set-TextValue "E42" "  -1.27%  "
Set-TextValue "D43" "2.18"
Set-TextValue "E43" "  -3.39%  "
Set-TextValue "D44" "0.926"
Set-TextValue "E44" "  -1.90%  "
Set-TextValue "D45" "46.59"
Set-TextValue "E45" "  -2.24%  "
Set-TextValue "D46" "2.72"
Set-TextValue "E46" "  -0.10%  "
Set-TextValue "E47" "  -3.13%  "
Set-TextValue "E48" "  -3.72%  "
Set-TextValue "E49" "  +0.20%  "
Set-TextValue "D50" "26.76"
Set-TextValue "E50" "  -5.17%  "
Set-TextValue "D51" "1.06"
Set-TextValue "E51" "  -5.62%  "
